# Updated cryptos list on Fri Jun  9 06:37:09 UTC 2023 with GitHub Actions
# Refreshes Price (D) / Volume(1h) (E) figures for each coin row, and swaps
# the ARBITRUM/ImmutableX (rows 34-35) and Cronos/Algorand (rows 48-49) pairs
# back into their updated ranking order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.547.91'
$ws.Range('E2').Value = '  +0.53%  '
$ws.Range('D3').Value = '1.841.53'
$ws.Range('E3').Value = '  +0.06%  '
$ws.Range('E4').Value = '  +0.04%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '259.77'
$c.Style = "Normal"
$ws.Range('E5').Value = '  -0.05%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '1.001'
$c.Style = "Normal"
$ws.Range('E6').Value = '  -0.03%  '
$ws.Range('E7').Value = '  +0.91%  '
$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '0.3180'
$c.Style = "Normal"
$ws.Range('E8').Value = '  -2.23%  '
$ws.Range('E9').Value = '  +0.40%  '
$ws.Range('E10').Value = '  +1.12%  '
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '0.7817'
$c.Style = "Normal"
$ws.Range('E11').Value = '  +2.40%  '
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '0.07777'
$c.Style = "Normal"
$ws.Range('D13').Value = '1.809.20'
$ws.Range('E13').Value = '  -1.64%  '
$ws.Range('E14').Value = '  -0.09%  '
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '5.019'
$c.Style = "Normal"
$ws.Range('E15').Value = '  +0.04%  '
$ws.Range('E16').Value = '  -0.04%  '
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '13.87'
$c.Style = "Normal"
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '1.000'
$c.Style = "Normal"
$ws.Range('E18').Value = '  -0.10%  '
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '0.000007944'
$c.Style = "Normal"
$ws.Range('E19').Value = '  +0.10%  '
$ws.Range('D20').Value = '26.574.69'
$ws.Range('E20').Value = '  +0.55%  '
$ws.Range('D21').Value = '2.069.68'
$ws.Range('E21').Value = '  -0.08%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '4.612'
$c.Style = "Normal"
$ws.Range('E22').Value = '  +1.25%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '5.977'
$c.Style = "Normal"
$ws.Range('E23').Value = '  +0.57%  '
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '9.342'
$c.Style = "Normal"
$ws.Range('E24').Value = '  -1.10%  '
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '142.72'
$c.Style = "Normal"
$ws.Range('E25').Value = '  -1.16%  '
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '2.216'
$c.Style = "Normal"
$ws.Range('E26').Value = '  +0.42%  '
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '1.670'
$c.Style = "Normal"
$ws.Range('E27').Value = '  +1.59%  '
$ws.Range('E28').Value = '  -0.28%  '
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '111.69'
$c.Style = "Normal"
$ws.Range('E29').Value = '  +0.64%  '
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '4.190'
$c.Style = "Normal"
$ws.Range('E30').Value = '  +0.71%  '
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '0.08728'
$c.Style = "Normal"
$ws.Range('E31').Value = '  +0.30%  '
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '4.077'
$c.Style = "Normal"
$ws.Range('E32').Value = '  -0.89%  '
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '0.04890'
$c.Style = "Normal"
$ws.Range('E33').Value = '  +2.68%  '
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '0.7250'
$c.Style = "Normal"
$ws.Range('E34').Value = '  +3.46%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '1.133'
$c.Style = "Normal"
$ws.Range('E35').Value = '  +1.08%  '
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '2.857'
$c.Style = "Normal"
$ws.Range('E36').Value = '  +0.55%  '
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '3.093'
$c.Style = "Normal"
$ws.Range('E37').Value = '  +1.28%  '
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '2.230'
$c.Style = "Normal"
$ws.Range('E38').Value = '  +2.83%  '
$ws.Range('E39').Value = '  -0.30%  '
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '0.4816'
$c.Style = "Normal"
$ws.Range('E40').Value = '  +0.28%  '
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '0.8977'
$c.Style = "Normal"
$ws.Range('E41').Value = '  +0.92%  '
$ws.Range('E42').Value = '  -0.50%  '
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '5.920'
$c.Style = "Normal"
$ws.Range('E43').Value = '  -2.30%  '
$ws.Range('E44').Value = '  -0.01%  '
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '7.657'
$c.Style = "Normal"
$ws.Range('E45').Value = '  +0.40%  '
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '0.4169'
$c.Style = "Normal"
$ws.Range('E46').Value = '  +1.68%  '
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '8.988'
$c.Style = "Normal"
$ws.Range('E47').Value = '  +0.57%  '
$ws.Range('B48').Value = 'Algorand'
$ws.Range('C48').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '0.1233'
$c.Style = "Normal"
$ws.Range('E48').Value = '  +1.51%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '0.05837'
$c.Style = "Normal"
$ws.Range('E49').Value = '  -0.47%  '
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '34.90'
$c.Style = "Normal"
$ws.Range('E50').Value = '  -0.26%  '
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '0.8932'
$c.Style = "Normal"
$ws.Range('E51').Value = '  +1.07%  '
